{"js": "const replacements = [\n  [\"205\u00d72=410\", \"823\u00d74=3292\"],\n  [\"486\u00d74=1944\", \"856\u00d74=3424\"],\n  [\"573\u00d78=4584\", \"451\u00d78=3608\"],\n  [\"495\u00d77=3465\", \"167\u00d79=1503\"],\n  [\"651\u00d72=1302\", \"186\u00d73=558\"],\n  [\"821\u00d74=3284\", \"517\u00d73=1551\"],\n  [\"801\u00d79=7209\", \"419\u00d73=1257\"],\n  [\"442\u00d77=3094\", \"546\u00d78=4368\"],\n  [\"309\u00d74=1236\", \"841\u00d77=5887\"],\n  [\"855\u00d72=1710\", \"713\u00d72=1426\"],\n  [\"456\u00d72=912\", \"878\u00d78=7024\"],\n  [\"927\u00d72=1854\", \"822\u00d79=7398\"],\n  [\"492\u00d78=3936\", \"454\u00d73=1362\"],\n  [\"231\u00d73=693\", \"771\u00d75=3855\"],\n  [\"804\u00d79=7236\", \"841\u00d74=3364\"],\n  [\"282\u00d75=1410\", \"992\u00d76=5952\"],\n  [\"254\u00d74=1016\", \"648\u00d78=5184\"],\n  [\"900\u00d75=4500\", \"290\u00d77=2030\"],\n  [\"105\u00d78=840\", \"259\u00d75=1295\"],\n  [\"620\u00d74=2480\", \"106\u00d79=954\"],\n  [\"295\u00d79=2655\", \"665\u00d75=3325\"],\n  [\"230\u00d76=1380\", \"909\u00d72=1818\"],\n  [\"369\u00d77=2583\", \"528\u00d76=3168\"],\n  [\"311\u00d79=2799\", \"517\u00d78=4136\"],\n  [\"614\u00d77=4298\", \"406\u00d73=1218\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"205\u00d72=410\", \"823\u00d74=3292\"),\n    @(\"486\u00d74=1944\", \"856\u00d74=3424\"),\n    @(\"573\u00d78=4584\", \"451\u00d78=3608\"),\n    @(\"495\u00d77=3465\", \"167\u00d79=1503\"),\n    @(\"651\u00d72=1302\", \"186\u00d73=558\"),\n    @(\"821\u00d74=3284\", \"517\u00d73=1551\"),\n    @(\"801\u00d79=7209\", \"419\u00d73=1257\"),\n    @(\"442\u00d77=3094\", \"546\u00d78=4368\"),\n    @(\"309\u00d74=1236\", \"841\u00d77=5887\"),\n    @(\"855\u00d72=1710\", \"713\u00d72=1426\"),\n    @(\"456\u00d72=912\", \"878\u00d78=7024\"),\n    @(\"927\u00d72=1854\", \"822\u00d79=7398\"),\n    @(\"492\u00d78=3936\", \"454\u00d73=1362\"),\n    @(\"231\u00d73=693\", \"771\u00d75=3855\"),\n    @(\"804\u00d79=7236\", \"841\u00d74=3364\"),\n    @(\"282\u00d75=1410\", \"992\u00d76=5952\"),\n    @(\"254\u00d74=1016\", \"648\u00d78=5184\"),\n    @(\"900\u00d75=4500\", \"290\u00d77=2030\"),\n    @(\"105\u00d78=840\", \"259\u00d75=1295\"),\n    @(\"620\u00d74=2480\", \"106\u00d79=954\"),\n    @(\"295\u00d79=2655\", \"665\u00d75=3325\"),\n    @(\"230\u00d76=1380\", \"909\u00d72=1818\"),\n    @(\"369\u00d77=2583\", \"528\u00d76=3168\"),\n    @(\"311\u00d79=2799\", \"517\u00d78=4136\"),\n    @(\"614\u00d77=4298\", \"406\u00d73=1218\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Not found: $old\"\n    }\n}\n"}
